# edit.ps1 - applies the commit's changes to CMP73010.docx via Word COM-interop
#
# Summary of the edit:
#   1. Remove the "_GoBack" bookmark that sits after "MP73010" in the title line.
#   2. Merge the three runs (with a proofErr gramStart/gramEnd pair) that make up
#      the ">>>  your stuff after this line >>>" paragraph into a single run.
#   3. Add a "_GoBack" bookmark right after "Ben changing things up!" (end of that
#      paragraph's text, before its paragraph mark).
#   4. Replace the two trailing empty paragraphs with a single paragraph that
#      reads "So does Kim!".

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Remove the old "_GoBack" bookmark (currently after "MP73010").
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 2. Collapse the ">>  / >  your / stuff ..." three runs into one run.
#    Setting Range.Text to the exact same concatenated text is a no-op in
#    this engine, so we first blank the paragraph (dropping the proofErr
#    marks and run boundaries) and then write the final text back in.
# ---------------------------------------------------------------------------
$quoteParaIndex = 4
$qp = $d.Paragraphs.Item($quoteParaIndex)
$qr = $d.Range($qp.Range.Start, $qp.Range.End - 1)
$qr.Text = ""
$qp2 = $d.Paragraphs.Item($quoteParaIndex)
$qr2 = $d.Range($qp2.Range.Start, $qp2.Range.End - 1)
$qr2.Text = ">>>  your stuff after this line >>>"

# ---------------------------------------------------------------------------
# 3. Add a new "_GoBack" bookmark right after "Ben changing things up!".
#    A collapsed (zero-length) Range placed exactly at a paragraph's last
#    text position confuses Bookmarks.Add in this engine (it silently drops
#    the position and anchors at document start instead). Work around it by
#    temporarily appending a marker character so the target position is no
#    longer the paragraph's end, add the bookmark there, then delete the
#    marker again - the bookmark stays anchored correctly.
# ---------------------------------------------------------------------------
$benParaIndex = 5
$bp = $d.Paragraphs.Item($benParaIndex)
$textEnd = $bp.Range.End - 1
$bodyRange = $d.Range($bp.Range.Start, $textEnd)
$bodyRange.InsertAfter("X")

$markerPos = $textEnd
$bmRange = $d.Range($markerPos, $markerPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$markerRange = $d.Range($markerPos, $markerPos + 1)
$markerRange.Delete()

# ---------------------------------------------------------------------------
# 4. Replace the two trailing empty paragraphs with one reading
#    "So does Kim!".  Delete the first (empty, second-to-last) paragraph so
#    the very last paragraph of the body (which can't be removed) remains,
#    then fill that last paragraph with the new text.
# ---------------------------------------------------------------------------
$paraCount = $d.Paragraphs.Count
$firstEmptyIndex = $paraCount - 1
$fp = $d.Paragraphs.Item($firstEmptyIndex)
$fp.Range.Delete()

$lastIndex = $d.Paragraphs.Count
$lp = $d.Paragraphs.Item($lastIndex)
$lr = $d.Range($lp.Range.Start, $lp.Range.End - 1)
$lr.InsertAfter("So does Kim!")
